# Applies the commit "updated project doc again":
#   1. Removes the old _GoBack bookmark that sat after
#      "Some questions we hope to answer"
#   2. Splits the indentation/content of the "Data Sets" section:
#        - The "Stock Prices ..." paragraph becomes a hanging/left indented
#          paragraph and gains a new trailing sentence about approval ratings.
#        - The "We plan to pull ..." paragraph keeps its first-line indent,
#          and a new _GoBack bookmark is dropped in the middle of "Kaggle"
#          (Word's "last edit" marker ends up there after the edit).

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark -----------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Re-indent the "Stock Prices" paragraph & append new sentence ------
$stockPara = $d.Paragraphs(15)
$stockPara.Range.Find.Execute("Stock Prices") | Out-Null
$stockPara.Format.LeftIndent = 36
$stockPara.Format.FirstLineIndent = 0

# Append the new trailing run text just before the paragraph mark.
$stockEnd = $stockPara.Range
$insertPoint = $d.Range($stockEnd.End - 1, $stockEnd.End - 1)
$insertPoint.InsertAfter(" " + [char]8211 + " approval ratings for the President, the house & senate")

# --- 3. Insert a fresh _GoBack bookmark inside "Kaggle" -> "Kag|gle" ------
$pullPara = $d.Paragraphs(16)
$pullRange = $pullPara.Range.Duplicate
$pullRange.Find.Execute("Kaggle") | Out-Null
$splitAt = $pullRange.Start + 3
$bmRange = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $bmRange)
